$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression (name unchanged, B2 unchanged, C2/D2 slightly adjusted)
$ws.Range("C2").Value = 0.3334337871158492
$ws.Range("D2").Value = 0.3334337871158494

# Row 3 - RandomForestRegressor (name unchanged, values changed)
$ws.Range("B3").Value = 0.02521916229754797
$ws.Range("C3").Value = 0.02501915099426874
$ws.Range("D3").Value = 0.0948544579762262

# Row 4 - renamed GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.03052248170867583
$ws.Range("C4").Value = 0.03067768909371772
$ws.Range("D4").Value = 0.06089497414405966

# Row 5 - renamed AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.01557545462918059
$ws.Range("C5").Value = 0.01665501732200664
$ws.Range("D5").Value = 0.0208402466659823
